# Update vm_pu.xlsx results: case with 380 kV done
# Sets the bus-voltage (vm_pu) values for rows 2-25, columns B-F and I-M
# to the values produced by the re-run power-flow case (slack bus voltage
# changed from 1.05 p.u. to 1.02 p.u. in column B, cascading through the
# rest of the result columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.037478228321651
$ws.Cells.Item(2, 4).Value = 1.038879757060748
$ws.Cells.Item(2, 5).Value = 1.036232274417504
$ws.Cells.Item(2, 6).Value = 1.036181704921944
$ws.Cells.Item(2, 9).Value = 1.037523701088701
$ws.Cells.Item(2, 10).Value = 1.04258088169361
$ws.Cells.Item(2, 11).Value = 1.041666532877474
$ws.Cells.Item(2, 12).Value = 1.039026604927106
$ws.Cells.Item(2, 13).Value = 1.038976180141608
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038703356110267
$ws.Cells.Item(3, 4).Value = 1.039788269039322
$ws.Cells.Item(3, 5).Value = 1.037282930358146
$ws.Cells.Item(3, 6).Value = 1.038037252531344
$ws.Cells.Item(3, 9).Value = 1.037867201169698
$ws.Cells.Item(3, 10).Value = 1.043449020560751
$ws.Cells.Item(3, 11).Value = 1.04238488913646
$ws.Cells.Item(3, 12).Value = 1.039886180046278
$ws.Cells.Item(3, 13).Value = 1.040638502504299
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.039495071171303
$ws.Cells.Item(4, 4).Value = 1.040375253318094
$ws.Cells.Item(4, 5).Value = 1.037962170088302
$ws.Cells.Item(4, 6).Value = 1.03923668466756
$ws.Cells.Item(4, 9).Value = 1.038087786329221
$ws.Cells.Item(4, 10).Value = 1.044009249485433
$ws.Cells.Item(4, 11).Value = 1.042848203570965
$ws.Cells.Item(4, 12).Value = 1.040441182765688
$ws.Cells.Item(4, 13).Value = 1.041712491584502
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03982766719519
$ws.Cells.Item(5, 4).Value = 1.040621812321287
$ws.Cells.Item(5, 5).Value = 1.038247580570004
$ws.Cells.Item(5, 6).Value = 1.039740642732961
$ws.Cells.Item(5, 9).Value = 1.038180118931455
$ws.Cells.Item(5, 10).Value = 1.044244410191316
$ws.Cells.Item(5, 11).Value = 1.043042621850265
$ws.Cells.Item(5, 12).Value = 1.040674221462992
$ws.Cells.Item(5, 13).Value = 1.042163613164486
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.039883497523653
$ws.Cells.Item(6, 4).Value = 1.040663198464676
$ws.Cells.Item(6, 5).Value = 1.038295494009012
$ws.Cells.Item(6, 6).Value = 1.039825243237108
$ws.Cells.Item(6, 9).Value = 1.038195598482827
$ws.Cells.Item(6, 10).Value = 1.044283873716584
$ws.Cells.Item(6, 11).Value = 1.043075244546273
$ws.Cells.Item(6, 12).Value = 1.04071333308673
$ws.Cells.Item(6, 13).Value = 1.042239336276636
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.039499516277262
$ws.Cells.Item(7, 4).Value = 1.040378548670602
$ws.Cells.Item(7, 5).Value = 1.037965984311102
$ws.Cells.Item(7, 6).Value = 1.039243419678419
$ws.Cells.Item(7, 9).Value = 1.038089021657518
$ws.Cells.Item(7, 10).Value = 1.04401239312313
$ws.Cells.Item(7, 11).Value = 1.042850802805613
$ws.Cells.Item(7, 12).Value = 1.040444297754391
$ws.Cells.Item(7, 13).Value = 1.041718520987033
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037892480720221
$ws.Cells.Item(8, 4).Value = 1.039186977354732
$ws.Cells.Item(8, 5).Value = 1.036587475526745
$ws.Cells.Item(8, 6).Value = 1.036809057834881
$ws.Cells.Item(8, 9).Value = 1.037640137716134
$ws.Cells.Item(8, 10).Value = 1.042874588877741
$ws.Cells.Item(8, 11).Value = 1.041909618871735
$ws.Cells.Item(8, 12).Value = 1.03931735249519
$ws.Cells.Item(8, 13).Value = 1.039538316517614
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.035052638911062
$ws.Cells.Item(9, 4).Value = 1.037080400641472
$ws.Cells.Item(9, 5).Value = 1.034153602216399
$ws.Cells.Item(9, 6).Value = 1.032509452394138
$ws.Cells.Item(9, 9).Value = 1.036836200442944
$ws.Cells.Item(9, 10).Value = 1.040857887343311
$ws.Cells.Item(9, 11).Value = 1.040239458036682
$ws.Cells.Item(9, 12).Value = 1.037322214968798
$ws.Cells.Item(9, 13).Value = 1.035683458989513
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033153732998444
$ws.Cells.Item(10, 4).Value = 1.035671240021467
$ws.Cells.Item(10, 5).Value = 1.03252762947213
$ws.Cells.Item(10, 6).Value = 1.029635632796225
$ws.Cells.Item(10, 9).Value = 1.036291450377111
$ws.Cells.Item(10, 10).Value = 1.039505321742256
$ws.Cells.Item(10, 11).Value = 1.039118013627882
$ws.Cells.Item(10, 12).Value = 1.035985680465977
$ws.Cells.Item(10, 13).Value = 1.033104125688684
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032330077603943
$ws.Cells.Item(11, 4).Value = 1.035059890107805
$ws.Cells.Item(11, 5).Value = 1.031822718662249
$ws.Cells.Item(11, 6).Value = 1.028389316112926
$ws.Cells.Item(11, 9).Value = 1.03605346236248
$ws.Cells.Item(11, 10).Value = 1.038917679694671
$ws.Cells.Item(11, 11).Value = 1.038630482309412
$ws.Cells.Item(11, 12).Value = 1.035405377490515
$ws.Cells.Item(11, 13).Value = 1.031984864960498
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032023917101108
$ws.Cells.Item(12, 4).Value = 1.03483262802729
$ws.Cells.Item(12, 5).Value = 1.031560751360062
$ws.Cells.Item(12, 6).Value = 1.027926074054699
$ws.Cells.Item(12, 9).Value = 1.035964744588052
$ws.Cells.Item(12, 10).Value = 1.038699102958051
$ws.Cells.Item(12, 11).Value = 1.038449097120726
$ws.Cells.Item(12, 12).Value = 1.035189586953487
$ws.Cells.Item(12, 13).Value = 1.031568749739562
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03208959952994
$ws.Cells.Item(13, 4).Value = 1.034881384688183
$ws.Cells.Item(13, 5).Value = 1.031616950271464
$ws.Cells.Item(13, 6).Value = 1.02802545505639
$ws.Cells.Item(13, 9).Value = 1.035983789290177
$ws.Cells.Item(13, 10).Value = 1.038746002064726
$ws.Cells.Item(13, 11).Value = 1.038488018229669
$ws.Cells.Item(13, 12).Value = 1.035235885682573
$ws.Cells.Item(13, 13).Value = 1.031658024900671
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032304774764128
$ws.Cells.Item(14, 4).Value = 1.035041108233244
$ws.Cells.Item(14, 5).Value = 1.031801067078779
$ws.Cells.Item(14, 6).Value = 1.028351030674069
$ws.Cells.Item(14, 9).Value = 1.036046135429787
$ws.Cells.Item(14, 10).Value = 1.038899618223115
$ws.Cells.Item(14, 11).Value = 1.038615494975892
$ws.Cells.Item(14, 12).Value = 1.035387545097773
$ws.Cells.Item(14, 13).Value = 1.031950476407016
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032437322249959
$ws.Cells.Item(15, 4).Value = 1.035139495301461
$ws.Cells.Item(15, 5).Value = 1.031914489910701
$ws.Cells.Item(15, 6).Value = 1.028551587780544
$ws.Cells.Item(15, 9).Value = 1.036084506695577
$ws.Cells.Item(15, 10).Value = 1.038994226292944
$ws.Cells.Item(15, 11).Value = 1.038693998508673
$ws.Cells.Item(15, 12).Value = 1.035480955535515
$ws.Cells.Item(15, 13).Value = 1.032130615768713
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033208365989241
$ws.Cells.Item(16, 4).Value = 1.035711788300682
$ws.Cells.Item(16, 5).Value = 1.032574393770988
$ws.Cells.Item(16, 6).Value = 1.029718304792633
$ws.Cells.Item(16, 9).Value = 1.036307200300687
$ws.Cells.Item(16, 10).Value = 1.039544279715603
$ws.Cells.Item(16, 11).Value = 1.039150328369005
$ws.Cells.Item(16, 12).Value = 1.036024159780079
$ws.Cells.Item(16, 13).Value = 1.033178355865571
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.033691638300562
$ws.Cells.Item(17, 4).Value = 1.036070456061115
$ws.Cells.Item(17, 5).Value = 1.032988102779286
$ws.Cells.Item(17, 6).Value = 1.030449627949475
$ws.Cells.Item(17, 9).Value = 1.036446324588335
$ws.Cells.Item(17, 10).Value = 1.039888782855147
$ws.Cells.Item(17, 11).Value = 1.039436050856392
$ws.Cells.Item(17, 12).Value = 1.036364473330154
$ws.Cells.Item(17, 13).Value = 1.033834926227527
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033973386548659
$ws.Cells.Item(18, 4).Value = 1.036279547829997
$ws.Cells.Item(18, 5).Value = 1.033229330113582
$ws.Cells.Item(18, 6).Value = 1.030876010739298
$ws.Cells.Item(18, 9).Value = 1.036527270196175
$ws.Cells.Item(18, 10).Value = 1.040089535691118
$ws.Cells.Item(18, 11).Value = 1.039602521125658
$ws.Cells.Item(18, 12).Value = 1.03656282077397
$ws.Cells.Item(18, 13).Value = 1.03421766263167
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034069432460876
$ws.Cells.Item(19, 4).Value = 1.036350823634859
$ws.Cells.Item(19, 5).Value = 1.033311568582599
$ws.Cells.Item(19, 6).Value = 1.031021365187221
$ws.Cells.Item(19, 9).Value = 1.03655483613783
$ws.Cells.Item(19, 10).Value = 1.040157955110677
$ws.Cells.Item(19, 11).Value = 1.03965925159046
$ws.Cells.Item(19, 12).Value = 1.036630426504745
$ws.Cells.Item(19, 13).Value = 1.034348127223342
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033639801895292
$ws.Cells.Item(20, 4).Value = 1.036031986123459
$ws.Cells.Item(20, 5).Value = 1.032943724241492
$ws.Cells.Item(20, 6).Value = 1.03037118320093
$ws.Cells.Item(20, 9).Value = 1.036431418903939
$ws.Cells.Item(20, 10).Value = 1.039851840613081
$ws.Cells.Item(20, 11).Value = 1.039405414885035
$ws.Cells.Item(20, 12).Value = 1.036327976624191
$ws.Cells.Item(20, 13).Value = 1.033764506322249
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032241417088069
$ws.Cells.Item(21, 4).Value = 1.034994078610522
$ws.Cells.Item(21, 5).Value = 1.031746852943334
$ws.Cells.Item(21, 6).Value = 1.028255165289579
$ws.Cells.Item(21, 9).Value = 1.036027784850193
$ws.Cells.Item(21, 10).Value = 1.038854390415821
$ws.Cells.Item(21, 11).Value = 1.038577964402045
$ws.Cells.Item(21, 12).Value = 1.035342891838983
$ws.Cells.Item(21, 13).Value = 1.031864367096296
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031360931828848
$ws.Cells.Item(22, 4).Value = 1.03434046485372
$ws.Cells.Item(22, 5).Value = 1.030993566619421
$ws.Cells.Item(22, 6).Value = 1.026922971792768
$ws.Cells.Item(22, 9).Value = 1.03577216087322
$ws.Cells.Item(22, 10).Value = 1.038225514343719
$ws.Cells.Item(22, 11).Value = 1.038056008862621
$ws.Cells.Item(22, 12).Value = 1.034722138562486
$ws.Cells.Item(22, 13).Value = 1.030667515879548
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.031827815462166
$ws.Cells.Item(23, 4).Value = 1.034687057510168
$ws.Cells.Item(23, 5).Value = 1.031392971807026
$ws.Cells.Item(23, 6).Value = 1.027629365163806
$ws.Cells.Item(23, 9).Value = 1.035907847294798
$ws.Cells.Item(23, 10).Value = 1.038559059685092
$ws.Cells.Item(23, 11).Value = 1.03833286999086
$ws.Cells.Item(23, 12).Value = 1.03505134471237
$ws.Cells.Item(23, 13).Value = 1.031302198139536
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.0336632249564
$ws.Cells.Item(24, 4).Value = 1.036049369381574
$ws.Cells.Item(24, 5).Value = 1.032963777246038
$ws.Cells.Item(24, 6).Value = 1.03040662957843
$ws.Cells.Item(24, 9).Value = 1.036438154768965
$ws.Cells.Item(24, 10).Value = 1.039868533808963
$ws.Cells.Item(24, 11).Value = 1.039419258538845
$ws.Cells.Item(24, 12).Value = 1.036344468383116
$ws.Cells.Item(24, 13).Value = 1.03379632675567
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035787786781522
$ws.Cells.Item(25, 4).Value = 1.037625831710977
$ws.Cells.Item(25, 5).Value = 1.034783400594856
$ws.Cells.Item(25, 6).Value = 1.033622256976015
$ws.Cells.Item(25, 9).Value = 1.03704558066865
$ws.Cells.Item(25, 10).Value = 1.04138066548835
$ws.Cells.Item(25, 11).Value = 1.040672633958485
$ws.Cells.Item(25, 12).Value = 1.037839128192535
$ws.Cells.Item(25, 13).Value = 1.03668164605183
